$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.690.79"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.14%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.901.87"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.84%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.36%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.39"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.33%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5188"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +7.72%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3781"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.17%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07244"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.07%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.10"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.61%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8955"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.32%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07639"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.64%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.890.21"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.435"
$ws.Range("D14").ClearFormats()

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.98"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.43%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.40%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008713"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.86%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9993"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.33%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "27.726.69"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.01%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.46"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.130"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.34%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.113.63"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.81"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.13%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.573"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.47"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.10%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.864"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.13%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.163"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.29%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.29"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.29%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.61"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.04%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.836"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.10%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08960"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.43%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.174"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.70%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.238"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.48%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.802"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.99%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7751"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.51%  "

$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02081"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.52%  "

$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.600"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.75%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.059"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.89%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.088"
$ws.Range("D39").ClearFormats()

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5491"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.41%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.05283"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.78%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.649"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.19%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "113.06"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.76%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.448"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.83%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1505"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.54%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4784"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.50%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.41"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.90%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9996"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.35%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.610"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.05%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "66.44"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.27%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06003"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.92%  "
